$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Done" column (E) with a header and mark the three
# "DB connection" related tasks (SQL - DB rows for catalog/login) as done.
$ws.Range("E2").Value = "Done"
$ws.Range("E31").Value = "done"
$ws.Range("E32").Value = "done"
$ws.Range("E33").Value = "done"

# Update the view so the newly added column is visible, matching the
# scrolled/selected state recorded in the workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 2
$ws.Range("C32").Select()
